$d = $word.ActiveDocument

function Replace-All($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# --- "English" appears twice (hyperlink label + standalone run); same replacement both times ---
Replace-All "English" "ภาษาอังกฤษ"

# --- language list line ---
Replace-All " / Portuguese / French / Thai / Vietnamese / Spanish" " / ภาษาโปรตุเกส / ภาษาฝรั่งเศส /ภาษาไทย / ภาษาเวียดนาม / ภาษาสเปน"

# --- Brief table ---
Replace-All "Brief" "บทย่อ"
Replace-All "An email sent to partners in the target country whose documents failed our verification process. It will be sent via customer.io" "An email sent to partners in the target country whose documents failed our verification process. โดยมันจะถูกส่งผ่านทาง customer.io"
Replace-All "Target audience" "กลุ่มเป้าหมาย"

# --- Heading ---
Replace-All "Uh oh! Your documents couldn’t be verified" "โอ ไม่นะ! เอกสารของคุณไม่อาจผ่านการตรวจสอบยืนยันได้"

# --- Greeting ---
Replace-All "Hi " "สวัสดี "
Replace-All ", " " "

# --- Body ---
Replace-All "We regret to inform you that your documents have failed our verification process as we found the following issues with them: " "เราขออภัยที่ต้องแจ้งให้คุณทราบว่า เอกสารของคุณไม่ผ่านกระบวนการตรวจสอบยืนยันของเรา เนื่องจากเราพบปัญหาดังต่อไปนี้: "

# --- List item 1 ---
Replace-All "A copy of your vaccination certificate" "สำเนาใบรับรองการฉีดวัคซีนของคุณ"
Replace-All ": Document is unclear" ": เอกสารไม่ชัดเจน"

# --- Resubmission paragraph ---
Replace-All "Please resubmit the documents above by " "กรุณายื่นเอกสารข้างต้นอีกครั้งภายในวันที่ "
Replace-All " so we can proceed with the necessary arrangements." " เพื่อให้เราสามารถดำเนินการตามขั้นตอนที่จำเป็นได้"

# --- Contact paragraph (live chat / WhatsApp) ---
Replace-All "If you have any questions, please contact us via " "หากคุณมีคำถามใดๆ กรุณาติดต่อเราผ่านทาง "
Replace-All "live chat" "แชทสด"

# First " or " (between the live-chat and WhatsApp hyperlinks) -- replace only the inner
# word so the run keeps its own (non-hyperlink) formatting instead of bleeding in the
# neighbouring hyperlink's formatting.
$rOr1 = $d.Content
$foundOr1 = $rOr1.Find.Execute(" or ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundOr1) {
    $innerOr1 = $d.Range($rOr1.Start + 1, $rOr1.End - 1)
    $innerOr1.Find.Execute("or", $true, $false, $false, $false, $false, $true, 1, $false, "หรือทาง", 2) | Out-Null
}

Replace-All ". " " "

# --- Country manager paragraph ---
Replace-All "If you have any questions, please contact your country manager, " "หากคุณมีคำถามใดๆ โปรดติดต่อผู้จัดการประจำประเทศของคุณซึ่งได้แก่ "
Replace-All ", at " " ที่ "

# Second " or " (between [EMAIL ADDRESS] and [WHATSAPP NO])
$rOr2 = $d.Content
$foundOr2 = $rOr2.Find.Execute(" or ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundOr2) {
    $innerOr2 = $d.Range($rOr2.Start + 1, $rOr2.End - 1)
    $innerOr2.Find.Execute("or", $true, $false, $false, $false, $false, $true, 1, $false, "หรือ", 2) | Out-Null
}

Replace-All " (WhatsApp). " " (WhatsApp) "

# --- Comment text ---
$d.Comments(1).Range.Text = "เลือกอย่างใดอย่างหนึ่ง"
